# Apply the "Added scripts for HIV AIDS mortality after HAART" commit's
# actual data edits to the Psychosis medications workbook:
#   1) Fix a typo: "Pimvanserin" -> "Pimavanserin"
#   2) Add a new column G "FDA_approved_generic_or_branded_2024" to the
#      Psychosis_medications sheet, with Yes/No values per drug
#   3) Widen/adjust a couple of column widths on that sheet
#   4) Update the "Drug_class" Sources text on the Metadata sheet to
#      include the PubChem URL
#   5) Add a new metadata row describing the new FDA_approved column

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Psychosis_medications
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Psychosis_medications")

# Fix the "Pimvanserin" -> "Pimavanserin" typo (row 25, Treatment_name)
$ws.Range("A25").Value = "Pimavanserin"

# New column header
$ws.Range("G1").Value = "FDA_approved_generic_or_branded_2024"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Per-drug FDA approval (as of July 2024) status, row by row
$fdaApproved = @{
    2  = "Yes"  # Chlorpromazine
    3  = "No"   # Reserpine
    4  = "Yes"  # Perphenazine
    5  = "No"   # Triflupromazine
    6  = "Yes"  # Haloperidol
    7  = "Yes"  # Fluphenazine
    8  = "Yes"  # Thioridazine
    9  = "Yes"  # Thiothixene
    10 = "Yes"  # Pimozide
    11 = "Yes"  # Molindone
    12 = "Yes"  # Loxapine
    13 = "Yes"  # Clozapine
    14 = "Yes"  # Risperidone
    15 = "Yes"  # Olanzapine
    16 = "Yes"  # Quetiapine
    17 = "Yes"  # Ziprasidone
    18 = "Yes"  # Aripiprazole
    19 = "Yes"  # Paliperidone
    20 = "Yes"  # Iloperidone
    21 = "Yes"  # Asenapine
    22 = "Yes"  # Lurasidone
    23 = "Yes"  # Cariprazine
    24 = "Yes"  # Brexpiprazole
    25 = "Yes"  # Pimavanserin
    26 = "Yes"  # Lumateperone
}

for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 7).Value = $fdaApproved[$row]
    $ws.Cells.Item($row, 6).Copy()
    $ws.Cells.Item($row, 7).PasteSpecial(-4122)
}

# Column width tweaks: column C gets its own (slightly wider) width,
# separate from column D (they used to share one <col> span); column F
# narrows a touch to make room for the new column.
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(6).ColumnWidth = 23

# ---------------------------------------------------------------------
# Sheet 2: Metadata
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Update the Drug_class row's Sources cell to include the PubChem URL
$meta.Range("E4").Value = "PubChem MeSH Tree classification (https://pubchem.ncbi.nlm.nih.gov/)"

# Add a metadata row describing the new FDA_approved_generic_or_branded_2024 column
$meta.Range("A8").Value = "FDA_approved_generic_or_branded_2024"
$meta.Range("B8").Value = "Is the drug is still in use and approved by the FDA as of July 2024? (Note: includes any version, i.e. original brand drug or generic)"
$meta.Range("C8").Value = "category, e.g. Yes or No"
$meta.Range("D8").Value = "Yes"
$meta.Range("E8").Value = "Drugs@FDA (https://www.accessdata.fda.gov/scripts/cder/daf/)"

$meta.Range("A7:E7").Copy()
$meta.Range("A8:E8").PasteSpecial(-4122)
